$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update workbook title (shared string) to reflect November 2016 data
$ws.Range("A2").Value = "by End-Use Sector, Census Division, and State, November 2016"

# Update Relative Standard Error data values (table body)
$ws.Range("B4").Value = 0.22
$ws.Range("B5").Value = 0.22
$ws.Range("D5").Value = 5
$ws.Range("B6").Value = 0.28999999999999998
$ws.Range("B7").Value = 0.48
$ws.Range("D7").Value = 8
$ws.Range("B8").Value = 0.32
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 3
$ws.Range("F10").Value = 2
$ws.Range("B11").Value = 0.1
$ws.Range("C11").Value = 0.17
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.19
$ws.Range("B12").Value = 0.17
$ws.Range("C12").Value = 0.3
$ws.Range("F12").Value = 0.34
$ws.Range("B13").Value = 0.16
$ws.Range("C13").Value = 0.23
$ws.Range("D13").Value = 2
$ws.Range("F13").Value = 0.35
$ws.Range("B14").Value = 0.18
$ws.Range("C14").Value = 0.39
$ws.Range("F14").Value = 0.25
$ws.Range("B15").Value = 0.27
$ws.Range("F15").Value = 0.49
$ws.Range("B16").Value = 0.34
$ws.Range("D17").Value = 2
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 0.36
$ws.Range("D20").Value = 5
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 1
$ws.Range("B24").Value = 2
$ws.Range("D24").Value = 6
$ws.Range("F24").Value = 3
$ws.Range("C25").Value = 1
$ws.Range("B27").Value = 2
$ws.Range("B28").Value = 3
$ws.Range("B29").Value = 0.39
$ws.Range("C29").Value = 0.26
$ws.Range("F29").Value = 0.32
$ws.Range("D30").Value = 8
$ws.Range("B32").Value = 0.47
$ws.Range("C32").Value = 0.43
$ws.Range("D32").Value = 4
$ws.Range("F32").Value = 0.46
$ws.Range("B33").Value = 1
$ws.Range("B34").Value = 0.22
$ws.Range("C34").Value = 0.43
$ws.Range("D34").Value = 4
$ws.Range("F34").Value = 0.38
$ws.Range("B35").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("D36").Value = 2
$ws.Range("C37").Value = 0.34
$ws.Range("B38").Value = 0.14000000000000001
$ws.Range("D38").Value = 0.37
$ws.Range("F38").Value = 0.21
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = 1
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 1
$ws.Range("C44").Value = 0.4
$ws.Range("F44").Value = 0.47
$ws.Range("C45").Value = 1
$ws.Range("B46").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("F47").Value = 1
$ws.Range("C48").Value = 0.42
$ws.Range("C50").Value = 3
$ws.Range("D51").Value = 5
$ws.Range("C52").Value = 5
$ws.Range("D52").Value = 5
$ws.Range("C53").Value = 8
$ws.Range("C56").Value = 6
$ws.Range("C57").Value = 7
$ws.Range("F57").Value = 3
$ws.Range("B58").Value = 0.41
$ws.Range("B59").Value = 0.43
$ws.Range("D63").Value = 12
$ws.Range("B65").Value = 0.22
$ws.Range("C65").Value = 0.34
$ws.Range("F65").Value = 0.24
